# RP3_FLT_EFF_2024_Jan_Dec.xlsx update:
#  - Refresh the "data current as of" date on FLT_EFF_YY!B2 (the other
#    sheets pull it via =FLT_EFF_YY!B2 so they recalc automatically).
#  - Update a handful of efficiency figures on ERT_FLT_EFF_FAB (row 6-7).
#  - Backfill the previously-blank column D figures on ERT_FLT_EFF_LOC
#    for rows 6 through 33.

$wb = $excel.ActiveWorkbook

# --- FLT_EFF_YY: bump the reference date (drives the formula copies on
#     FLT_EFF_MM, ERT_FLT_EFF_FAB and ERT_FLT_EFF_LOC) ---
$wsYY = $wb.Worksheets.Item("FLT_EFF_YY")
$wsYY.Range("B2").Value = 45758

# --- ERT_FLT_EFF_FAB: corrected figures on rows 6 and 7 ---
$wsFAB = $wb.Worksheets.Item("ERT_FLT_EFF_FAB")
$wsFAB.Range("D6").Value = 0.024
$wsFAB.Range("B7").Value = 0.061
$wsFAB.Range("C7").Value = 0.0577
$wsFAB.Range("E7").Value = 0.0348

# --- ERT_FLT_EFF_LOC: fill in column D (rows 6-33), previously blank ---
$wsLOC = $wb.Worksheets.Item("ERT_FLT_EFF_LOC")
$wsLOC.Range("D6").Value = 0.0196
$wsLOC.Range("D7").Value = 0.03
$wsLOC.Range("D8").Value = 0.0225
$wsLOC.Range("D9").Value = 0.0146
$wsLOC.Range("D10").Value = 0.0384
$wsLOC.Range("D11").Value = 0.0205
$wsLOC.Range("D12").Value = 0.0114
$wsLOC.Range("D13").Value = 0.0122
$wsLOC.Range("D14").Value = 0.0088
$wsLOC.Range("D15").Value = 0.0283
$wsLOC.Range("D16").Value = 0.023
$wsLOC.Range("D17").Value = 0.0192
$wsLOC.Range("D18").Value = 0.0149
$wsLOC.Range("D19").Value = 0.0113
$wsLOC.Range("D20").Value = 0.0267
$wsLOC.Range("D21").Value = 0.0125
$wsLOC.Range("D22").Value = 0.0192
$wsLOC.Range("D23").Value = 0.018
$wsLOC.Range("D24").Value = 0.0262
$wsLOC.Range("D25").Value = 0.0155
$wsLOC.Range("D26").Value = 0.0165
$wsLOC.Range("D27").Value = 0.018
$wsLOC.Range("D28").Value = 0.0205
$wsLOC.Range("D29").Value = 0.0213
$wsLOC.Range("D30").Value = 0.0155
$wsLOC.Range("D31").Value = 0.0308
$wsLOC.Range("D32").Value = 0.0105
$wsLOC.Range("D33").Value = 0.0395
